# Update data: 2025-10-30 12:37
# Applies the refreshed market-health snapshot to the workbook:
#   1. Metadata!A2      - bump the "Last Updated" timestamp by one minute.
#   2. Top Gainers      - a new leader (MIDWESTLTD) enters the board at
#                         row 40, pushing every following row down by one;
#                         the row that falls off the bottom (KERNEX) is gone.
#   3. distance from Dma50 - refreshed "Distance From Sma50" readings
#                         (stock order unchanged).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet: bump the "Last Updated" timestamp.
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "30 Oct 2025, 12:37 PM"

# ---------------------------------------------------------------------
# 2. Top Gainers sheet: new row 40, rows 40-75 shift down to 41-76,
#    old row 76 drops off the bottom of the table.
# ---------------------------------------------------------------------
$gainers = $wb.Worksheets.Item("Top Gainers")

$gainersRows = @(
    @{Row=40; B="MIDWESTLTD"; C=3.6792; D=-0.9206; E="N/A"}
    @{Row=41; B="REDTAPE"; C=3.6536; D=3.553; E=-3.3779}
    @{Row=42; B="SHRINGARMS"; C=3.4862; D=4.6964; E=24.9205}
    @{Row=43; B="RSYSTEMS"; C=3.4508; D=4.612; E=6.9611}
    @{Row=44; B="BLS"; C=3.4321; D=0.3771; E=-0.8843}
    @{Row=45; B="ALICON"; C=3.4194; D=9.5502; E=14.9919}
    @{Row=46; B="CENTUM"; C=3.3549; D=4.2178; E=-1.3077}
    @{Row=47; B="GANESHCP"; C=3.3463; D=2.8184; E=2.3676}
    @{Row=48; B="OIL"; C=3.2937; D=3.5399; E=4.9662}
    @{Row=49; B="GMMPFAUDLR"; C=3.2614; D=7.7366; E=20.1095}
    @{Row=50; B="BGRENERGY"; C=3.24; D=-6.0691; E=74.87}
    @{Row=51; B="PFOCUS"; C=2.9919; D=0.2857; E=1.7392}
    @{Row=52; B="IVALUE"; C=2.9631; D=6.2091; E=-0.9857}
    @{Row=53; B="SPANDANA"; C=2.9101; D=4.4085; E=3.021}
    @{Row=54; B="NEULANDLAB"; C=2.8885; D=-1.4731; E=8.5313}
    @{Row=55; B="BPCL"; C=2.8871; D=8.3825; E=5.4468}
    @{Row=56; B="CARYSIL"; C=2.859; D=2.3423; E=11.247}
    @{Row=57; B="MFSL"; C=2.8312; D=2.8854; E=-0.8947000000000001}
    @{Row=58; B="JKLAKSHMI"; C=2.7853; D=4.7782; E=1.7903}
    @{Row=59; B="BEML"; C=2.7641; D=-0.1442; E=6.179}
    @{Row=60; B="NBCC"; C=2.7252; D=5.9719; E=10.5342}
    @{Row=61; B="IIFL"; C=2.7151; D=9.6774; E=18.8759}
    @{Row=62; B="POWERINDIA"; C=2.6999; D=7.0794; E=-0.2611}
    @{Row=63; B="OBEROIRLTY"; C=2.6672; D=3.4826; E=11.1877}
    @{Row=64; B="AHLUCONT"; C=2.6261; D=1.4385; E=-5.7576}
    @{Row=65; B="INOXGREEN"; C=2.548; D=10.5432; E=33.882}
    @{Row=66; B="ASHOKA"; C=2.5392; D=4.0551; E=6.6742}
    @{Row=67; B="VOLTAMP"; C=2.5248; D=2.404; E=2.1422}
    @{Row=68; B="REFEX"; C=2.5106; D=-0.055; E=1.9212}
    @{Row=69; B="DBCORP"; C=2.4854; D=5.1178; E=1.0771}
    @{Row=70; B="SDBL"; C=2.3903; D=0.9379999999999999; E=6.5186}
    @{Row=71; B="SUNDROP"; C=2.3659; D=2.2593; E=0.3706}
    @{Row=72; B="JKTYRE"; C=2.3577; D=5.389; E=21.3896}
    @{Row=73; B="BLUEDART"; C=2.29; D=20.7888; E=17.8766}
    @{Row=74; B="DATAMATICS"; C=2.2739; D=9.7554; E=18.3614}
    @{Row=75; B="GRAPHITE"; C=2.2679; D=15.8864; E=16.0951}
    @{Row=76; B="DBL"; C=2.2584; D=3.3396; E=4.4108}
)

foreach ($r in $gainersRows) {
    $gainers.Cells.Item($r.Row, 2).Value = $r.B
    $gainers.Cells.Item($r.Row, 3).Value = $r.C
    $gainers.Cells.Item($r.Row, 4).Value = $r.D
    $gainers.Cells.Item($r.Row, 5).Value = $r.E
}

# ---------------------------------------------------------------------
# 3. "distance from Dma50" sheet: refreshed readings (stock order is
#    unchanged, only column C values move).
# ---------------------------------------------------------------------
$dma50 = $wb.Worksheets.Item("distance from Dma50")

$dma50Rows = @(
    @{Row=2; C=10.054}
    @{Row=3; C=7.4818}
    @{Row=4; C=6.3286}
    @{Row=5; C=5.3439}
    @{Row=6; C=5.3248}
    @{Row=7; C=5.0368}
    @{Row=8; C=4.4295}
    @{Row=9; C=4.3445}
    @{Row=10; C=3.8775}
    @{Row=11; C=3.6409}
    @{Row=12; C=3.4167}
    @{Row=13; C=3.3733}
    @{Row=14; C=3.0593}
    @{Row=15; C=3.0339}
    @{Row=16; C=2.9435}
    @{Row=17; C=2.8144}
    @{Row=18; C=2.7176}
    @{Row=19; C=2.6922}
    @{Row=20; C=2.3349}
    @{Row=21; C=2.3075}
    @{Row=22; C=1.3682}
    @{Row=23; C=1.3139}
    @{Row=24; C=1.2864}
    @{Row=25; C=1.063}
    @{Row=26; C=0.9392}
    @{Row=28; C=0.525}
    @{Row=29; C=0.257}
    @{Row=30; C=-2.0498}
)

foreach ($r in $dma50Rows) {
    $dma50.Cells.Item($r.Row, 3).Value = $r.C
}
